# Apply calendar update: replace Trimester 1 course-specific rows with the
# new T1 orientation / weekly / census-date schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Date($cellRef, [int]$year, [int]$month, [int]$day) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "d-mmm"
    $cell.Value = (Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0)
}

# ---------------------------------------------------------------------
# Row 2: T1 Orientation
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "T1 Orientation"
$ws.Range("B2").Value = "T1 Orientation"
$ws.Range("C2").Value = $null
Set-Date "E2" 2025 9 8
$ws.Range("F2").ClearContents()
Set-Date "G2" 2025 9 14
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = $null
$ws.Range("L2").Value = $null
$ws.Range("D2").Value = "Lecture"
$ws.Range("M2").Value = "Transparent"

# ---------------------------------------------------------------------
# Row 3: T1 - Week 1
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "T1 - Week 1"
$ws.Range("B3").Value = "T1 - Week 1"
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = "Orientation"
Set-Date "E3" 2025 9 15
$ws.Range("F3").Value = $null
Set-Date "G3" 2025 9 21
$ws.Range("H3").Value = $null
$ws.Range("I3").ClearContents()
$ws.Range("J3").Value = $null
$ws.Range("K3").ClearContents()
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = "Transparent"

# ---------------------------------------------------------------------
# Row 4: T1 - Week 2
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "T1 - Week 2"
$ws.Range("B4").Value = "T1 - Week 2"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "Tutorial"
Set-Date "E4" 2025 9 22
$ws.Range("F4").ClearContents()
Set-Date "G4" 2025 9 28
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = $null
$ws.Range("M4").Value = "Transparent"

# ---------------------------------------------------------------------
# Row 5: T1 - Week 3
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "T1 - Week 3"
$ws.Range("B5").Value = "T1 - Week 3"
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "TEST"
Set-Date "E5" 2025 9 29
$ws.Range("F5").ClearContents()
Set-Date "G5" 2025 10 5
$ws.Range("I5").Value = $null
$ws.Range("J5").Value = $null
$ws.Range("M5").Value = "Transparent"

# ---------------------------------------------------------------------
# Row 6 & 7: Census Date rows (new rows). Populate the "Census Date"
# labels (column B) before the "T1 - ... Census Date" titles (column A)
# so that newly introduced shared strings are interned in the same
# order as the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "Census Date"
$ws.Range("B7").Value = "Census Date International"
$ws.Range("A7").Value = "T1 - International Census Date"
$ws.Range("A6").Value = "T1 - Domestic Census Date"

Set-Date "E6" 2025 9 28
$ws.Range("M6").Value = "Transparent"

Set-Date "E7" 2025 10 5
$ws.Range("M7").Value = "Transparent"

# ---------------------------------------------------------------------
# Final selection state, matching the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("F25").Select()
